$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 160
$ws1.Range("F4").Value = 576
$ws1.Range("F5").Value = 1804
$ws1.Range("F9").Value = 2234
$ws1.Range("F11").Value = 58
$ws1.Range("F12").Value = 153
$ws1.Range("F13").Value = 1384
$ws1.Range("F21").Value = 57
$ws1.Range("F23").Value = 52
$ws1.Range("F25").Value = 1400
$ws1.Range("F27").Value = 361

# --- Sheet "全部类型" (all types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 160
$ws4.Range("F4").Value = 576
$ws4.Range("F5").Value = 1804
$ws4.Range("F10").Value = 2234
$ws4.Range("F12").Value = 58
$ws4.Range("F13").Value = 153
$ws4.Range("F14").Value = 1384
$ws4.Range("F22").Value = 57
$ws4.Range("F24").Value = 52
$ws4.Range("F26").Value = 1400
$ws4.Range("F28").Value = 361
$ws4.Range("F29").Value = 0

$wb.Save()
